# digital_passport.xlsx -> oca-data-vault:0.2.3 example-data refresh
#
# 1. SAI (Self-Addressing Identifier) codes regenerated for several schema
#    fields: "Main" (G4/G7 share one field id, G12), "en" (E4, E7, E12) and
#    "fr" (E4, E7 -- E7 shares the same field id as "en"!E7 -- and E12).
# 2. Active sheet moves from "Main" to "fr", with a new selected cell on
#    each touched sheet.
# 3. A handful of rows (the ones whose wrapped SAI-code cell drives the row
#    height) settle on a new auto-fit height after the text refresh.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Main")
$wsEn   = $wb.Worksheets.Item("en")
$wsFr   = $wb.Worksheets.Item("fr")

# --- 1. Refreshed SAI codes -------------------------------------------------

# Main sheet: "utf-8" field id (used twice) and "nationality" field id
$wsMain.Range("G4").Value  = "SAI:EGyWgdQR9dW_I5oHlHBMoO9AA_eMeb2p3XzcCRCBbKCM"
$wsMain.Range("G7").Value  = "SAI:EGyWgdQR9dW_I5oHlHBMoO9AA_eMeb2p3XzcCRCBbKCM"
$wsMain.Range("G12").Value = "SAI:EmXip-eMwEeLWH9_AfH642_Zb-oGZXk7jI49qafz_XrU"

# "en" sheet: issuing-state(full), issuing-state(code), nationality
$wsEn.Range("E4").Value  = "SAI:Els6NxGvFfyL5aiBWR3j7YiaS7F4j4O-F0EIlZu-dO0g"
$wsEn.Range("E7").Value  = "SAI:EdxqlME_1Zt0Y_YJ3c0uMIzd41mSDATbH-rp7ElqhNes"
$wsEn.Range("E12").Value = "SAI:EXdSkFdYnAzZ2U2Qyo-q76CJMYelgV9NXN8GhmtY2ErY"

# "fr" sheet: issuing-state(full), issuing-state(code), nationality
$wsFr.Range("E4").Value  = "SAI:Emj736NIuEKdR-3sfXB4wfPokFzgE7uiRm6EXKaJRyE4"
$wsFr.Range("E7").Value  = "SAI:EdxqlME_1Zt0Y_YJ3c0uMIzd41mSDATbH-rp7ElqhNes"
$wsFr.Range("E12").Value = "SAI:Ez_5-oggNDr7gUFaz3GSof1y579gK0MKXIVesRjGzyqY"

# --- 2. Row-height settle after the text refresh ----------------------------

$wsMain.Rows.Item(4).RowHeight  = 23.85
$wsMain.Rows.Item(12).RowHeight = 23.85

$wsEn.Rows.Item(4).RowHeight  = 23.85
$wsEn.Rows.Item(12).RowHeight = 23.85

$wsFr.Rows.Item(4).RowHeight  = 23.85
$wsFr.Rows.Item(7).RowHeight  = 23.85
$wsFr.Rows.Item(12).RowHeight = 35.05

# --- 3. Selection / active-sheet changes ------------------------------------

$wsMain.Range("J18").Select()

$wsEn.Activate()
$wsEn.Range("E12").Select()

$wsFr.Activate()
$wsFr.Range("E6").Select()
